$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.165075540542603
$ws.Range("B1").Value = 2.510160684585571
$ws.Range("C1").Value = 6.726939678192139
$ws.Range("D1").Value = 2.052569389343262
$ws.Range("E1").Value = 1.212252378463745
